$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 129 ----
# Numeric columns: date, volume, high, low, open, close
$ws.Cells.Item(129, 1).Value = 45498.2916666667
$ws.Cells.Item(129, 2).Value = 0
$ws.Cells.Item(129, 3).Value = 1.87000000476837
$ws.Cells.Item(129, 4).Value = 1.87000000476837
$ws.Cells.Item(129, 5).Value = 1.87000000476837
$ws.Cells.Item(129, 6).Value = 1.87000000476837
$ws.Cells.Item(129, 8).Value = "KK.MI"

# Copy the date format (style) from the cell above, same as the rest of column A
$ws.Cells.Item(128, 1).Copy()
$ws.Cells.Item(129, 1).PasteSpecial(-4122)

# adj_close (G) is stored as a numeric-looking shared text string ("1.87000000476837")
# like the rest of the column; copy an existing identical text cell so it is written
# back as text (not re-parsed into a number) and without adding a new style.
$ws.Cells.Item(128, 7).Copy()
$ws.Cells.Item(129, 7).PasteSpecial()

# ---- Row 130 ----
$ws.Cells.Item(130, 1).Value = 45499.6031712963
$ws.Cells.Item(130, 2).Value = 3600
$ws.Cells.Item(130, 3).Value = 1.9099999666214
$ws.Cells.Item(130, 4).Value = 1.89999997615814
$ws.Cells.Item(130, 5).Value = 1.9099999666214
$ws.Cells.Item(130, 6).Value = 1.89999997615814
$ws.Cells.Item(130, 8).Value = "KK.MI"

$ws.Cells.Item(128, 1).Copy()
$ws.Cells.Item(130, 1).PasteSpecial(-4122)

# adj_close "1.89999997615814" already exists verbatim as text in row 28 (and others)
$ws.Cells.Item(28, 7).Copy()
$ws.Cells.Item(130, 7).PasteSpecial()

$excel.CutCopyMode = $false
